$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.415.94"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "1.842.89"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'239.40"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("D6").Value = "'0.6267"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D8").Value = "'0.07443"

$ws.Range("D9").Value = "'0.2896"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").Value = "'24.98"
$ws.Range("E10").Value = "  +2.10%  "

$ws.Range("D11").Value = "'0.07718"
$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").Value = "1.844.13"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").Value = "'4.976"
$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").Value = "'0.6765"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").Value = "'0.00001029"
$ws.Range("E15").Value = "  -2.52%  "

$ws.Range("D16").Value = "'81.81"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").Value = "'6.244"
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").Value = "29.454.45"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("D19").Value = "'232.94"
$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("D20").Value = "'12.33"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "'7.332"
$ws.Range("E22").Value = "  -2.03%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "'158.18"
$ws.Range("E24").Value = "  -0.82%  "

$ws.Range("D25").Value = "'8.502"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("E26").Value = "  -1.56%  "

$ws.Range("E27").Value = "  -1.06%  "

$ws.Range("D28").Value = "'0.07154"
$ws.Range("E28").Value = "  +10.76%  "

$ws.Range("D29").Value = "'1.471"
$ws.Range("E29").Value = "  +4.06%  "

$ws.Range("D30").Value = "'1.485"
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'4.048"
$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.039"
$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("D33").Value = "'1.822"
$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("D35").Value = "'0.6984"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Value = "'2.578"
$ws.Range("E36").Value = "  -0.23%  "

$ws.Range("D37").Value = "'6.952"
$ws.Range("E37").Value = "  +2.55%  "

$ws.Range("D38").Value = "'0.01843"
$ws.Range("E38").Value = "  +0.79%  "

$ws.Range("D39").Value = "'2.819"
$ws.Range("E39").Value = "  -0.77%  "

$ws.Range("D40").Value = "1.236.20"
$ws.Range("E40").Value = "  -2.50%  "

$ws.Range("D41").Value = "'0.9674"
$ws.Range("E41").Value = "  +6.29%  "

$ws.Range("D42").Value = "'1.000"

$ws.Range("D43").Value = "2.007.33"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "'101.01"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").Value = "'65.49"
$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "'1.725"
$ws.Range("E47").Value = "  -0.79%  "

$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("D49").Value = "'8.954"
$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("D50").Value = "'0.1137"
$ws.Range("E50").Value = "  -2.42%  "

$ws.Range("D51").Value = "'0.3904"
$ws.Range("E51").Value = "  -1.61%  "
